$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text-like numeric value in the D (Price) column without
# Excel auto-converting it to a floating point number / changing cell style.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "66.206.24"
$ws.Range("E2").Value = "  -1.68%  "
Set-TextValue "D3" "3.277.97"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "577.27"
$ws.Range("E5").Value = "  -0.22%  "
Set-TextValue "D6" "179.18"
$ws.Range("E6").Value = "  -2.34%  "
Set-TextValue "D7" "0.625"
$ws.Range("E7").Value = "  +3.45%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  -1.61%  "
Set-TextValue "D12" "3.847.27"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("E13").Value = "  -3.69%  "
Set-TextValue "D14" "66.246.26"
$ws.Range("E14").Value = "  -1.93%  "
Set-TextValue "D15" "26.37"
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "3.384.61"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0000163"
$ws.Range("E17").Value = "  -2.64%  "
Set-TextValue "D18" "436.15"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("E19").Value = "  -2.47%  "
Set-TextValue "D20" "13.07"
$ws.Range("E20").Value = "  -3.75%  "
$ws.Range("E21").Value = "  -4.25%  "
Set-TextValue "D22" "71.75"
$ws.Range("E22").Value = "  -3.07%  "
Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  -0.10%  "
Set-TextValue "D24" "3.418.92"
$ws.Range("E24").Value = "  -1.41%  "
Set-TextValue "D25" "0.504"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  +3.35%  "
$ws.Range("E27").Value = "  -4.97%  "
$ws.Range("E28").Value = "  -2.16%  "
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("E35").Value = "  -3.83%  "
Set-TextValue "D36" "157.61"
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("E37").Value = "  -5.21%  "
Set-TextValue "D38" "26.62"
$ws.Range("E38").Value = "  -1.61%  "
Set-TextValue "D39" "1.79"
$ws.Range("E39").Value = "  -3.05%  "
Set-TextValue "D40" "2.761.32"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("E41").Value = "  -2.42%  "
Set-TextValue "D42" "4.30"
$ws.Range("E42").Value = "  -3.51%  "
Set-TextValue "D43" "40.26"
$ws.Range("E43").Value = "  -0.22%  "
Set-TextValue "D44" "6.04"
$ws.Range("E44").Value = "  -3.04%  "
Set-TextValue "D45" "0.0656"
$ws.Range("E45").Value = "  -2.40%  "
Set-TextValue "D46" "321.32"
$ws.Range("E46").Value = "  -1.13%  "
Set-TextValue "D47" "2.29"
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("E48").Value = "  -5.74%  "
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("E50").Value = "  +2.08%  "
$ws.Range("E51").Value = "  -0.03%  "
